# The document carries two logos, each duplicated once for the primary
# header/footer and once for the "first page" header/footer:
#   - the BTEC logo (header1 / header2), currently tagged "image2.jpg"
#     and being renumbered to "image1.jpg"
#   - the Pearson logo (footer1 / footer2), currently tagged "image1.png"
#     and being renumbered to "image2.png"
#
# Word keeps that bookkeeping name (<wp:docPr>/<pic:cNvPr> @name) on the
# Shape side of the object model, not on InlineShape, so the idiomatic way
# to rename an inline picture without disturbing its inline layout is:
#   InlineShape -> ConvertToShape -> set .Name -> ConvertToInlineShape.

$d = $word.ActiveDocument
$sec = $d.Sections(1)

function Rename-InlinePicture($range, $newName) {
    $inline = $range.InlineShapes(1)
    $shape = $inline.ConvertToShape()
    $shape.Name = $newName
    [void]$shape.ConvertToInlineShape()
}

# Header 1 (primary) - BTec_Logo-Orange, id=1: image2.jpg -> image1.jpg
Rename-InlinePicture $sec.Headers(1).Range "image1.jpg"

# Header 2 (first page) - BTec_Logo-Orange, id=3: image2.jpg -> image1.jpg
Rename-InlinePicture $sec.Headers(2).Range "image1.jpg"

# Footer 1 (primary) - PearsonLogo, id=2: image1.png -> image2.png
Rename-InlinePicture $sec.Footers(1).Range "image2.png"

# Footer 2 (first page) - PearsonLogo, id=4: image1.png -> image2.png
Rename-InlinePicture $sec.Footers(2).Range "image2.png"
